$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2.71
$ws.Range("C2").Value = -21.65
$ws.Range("D2").Value = 813.7879881286622

# Row 3
$ws.Range("B3").Value = 1.72
$ws.Range("C3").Value = 7.61
$ws.Range("D3").Value = 1212.970706129145

# Row 4
$ws.Range("B4").Value = 1.67
$ws.Range("C4").Value = -38.2
$ws.Range("D4").Value = 255.225017937296

# Row 5
$ws.Range("B5").Value = -2.63
$ws.Range("C5").Value = -38.55
$ws.Range("D5").Value = 159.1260830577557

# Row 6
$ws.Range("B6").Value = 1.45
$ws.Range("C6").Value = 130.19
$ws.Range("D6").Value = 1308.354170588115

# Row 7
$ws.Range("B7").Value = 0.11
$ws.Range("C7").Value = -14.44
$ws.Range("D7").Value = 627.3380717086934

# Row 8
$ws.Range("B8").Value = 2.11
$ws.Range("C8").Value = -24.46
$ws.Range("D8").Value = 351.6657779184883

# Row 9
$ws.Range("B9").Value = 0.5600000000000001
$ws.Range("C9").Value = -66.56
$ws.Range("D9").Value = 403.6500091552734

# Row 10
$ws.Range("B10").Value = 0.79
$ws.Range("C10").Value = -4.76
$ws.Range("D10").Value = 576

# Row 11
$ws.Range("B11").Value = -0.25
$ws.Range("C11").Value = -5.88
$ws.Range("D11").Value = 478.3949788650513

# Row 12
$ws.Range("B12").Value = 9.609999999999999
$ws.Range("C12").Value = 42.77
$ws.Range("D12").Value = 848.2599136505127

# Row 13
$ws.Range("B13").Value = 3.65
$ws.Range("C13").Value = -25.62
$ws.Range("D13").Value = 1117.762165676077

# Row 14
$ws.Range("B14").Value = 2.63
$ws.Range("C14").Value = 48.21
$ws.Range("D14").Value = 1049.233612022287

# Row 15
$ws.Range("B15").Value = -2.31
$ws.Range("C15").Value = -36.76
$ws.Range("D15").Value = 444.2428796580075

# Row 16
$ws.Range("B16").Value = 3.27
$ws.Range("C16").Value = 46.72
$ws.Range("D16").Value = 654.261925610733

# Row 17
$ws.Range("B17").Value = 5.77
$ws.Range("C17").Value = -69.31999999999999
$ws.Range("D17").Value = 118.4223753226282

# Row 18
$ws.Range("B18").Value = 2.26
$ws.Range("C18").Value = 4.19
$ws.Range("D18").Value = 10418.69567542873

# Row 19 (Cash row; B19 stays "---")
$ws.Range("C19").Value = 1381.565621111907
$ws.Range("D19").Value = 970.8250862732211
